$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4698
$ws1.Range("F4").Value = 69
$ws1.Range("F6").Value = 240
$ws1.Range("F7").Value = 141
$ws1.Range("F8").Value = 217
$ws1.Range("F9").Value = 184
$ws1.Range("F10").Value = 1795
$ws1.Range("F11").Value = 328
$ws1.Range("F12").Value = 4121
$ws1.Range("F14").Value = 284

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 48

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4698
$ws4.Range("F5").Value = 69
$ws4.Range("F7").Value = 48
$ws4.Range("F8").Value = 240
$ws4.Range("F9").Value = 141
$ws4.Range("F10").Value = 217
$ws4.Range("F11").Value = 184
$ws4.Range("F14").Value = 1795
$ws4.Range("F15").Value = 328
$ws4.Range("F16").Value = 4121
$ws4.Range("F18").Value = 284
